# Update automatico via Actualizar 02-05-2021 17-55-00
#
# The monitoring sheet logs one 14-row "check" block (one row per monitored
# service) every time the checker script runs. This run:
#   - re-stamped the D column of the previous block (rows 702:715) with a
#     (microsecond-truncated) copy of its own timestamp, and
#   - appended a brand-new 14-row block (rows 716:729) with a fresh
#     timestamp, reusing the same Nombre/URL/Disponibilidad text for every
#     monitored service (same order as every previous block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rows 702:715 - D column gets rewritten with (effectively) the same
#    instant, losing a little sub-millisecond precision in the process.
# ---------------------------------------------------------------------
$refreshedTimestamp = 44232.72488630787
for ($r = 702; $r -le 715; $r++) {
    $ws.Range("D$r").Value = $refreshedTimestamp
}

# ---------------------------------------------------------------------
# 2) Append the new block: rows 716:729.
# ---------------------------------------------------------------------
$names = @(
    "Odoo", "Blackbox", "PowerBI", "Dropbox", "Odoo", "GEE",
    "UtilidadesOdoo", "Filtros Dashboard", "MapStore", "GeoServer",
    "Tomcat", "Shiny", "Github", "EZ Exporter"
)

# Text shown in column B (what ends up in the cell / shared string).
$displayUrls = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

# Hyperlink target (Address) - same as display text except MapStore, whose
# "#/" anchor is stored as the hyperlink's SubAddress instead.
$linkAddresses = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

$subAddresses = @("", "", "", "", "", "", "", "", "/", "", "", "", "", "")

$newTimestamp = 44232.74624696843
$firstNewRow = 716

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $firstNewRow + $i

    $ws.Range("A$row").Value = $names[$i]
    $ws.Range("B$row").Value = $displayUrls[$i]
    $ws.Range("C$row").Value = "Disponible"

    $ws.Range("D$row").Value = $newTimestamp
    $ws.Range("D$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $cell = $ws.Range("B$row")
    $ws.Hyperlinks.Add($cell, $linkAddresses[$i], $subAddresses[$i]) | Out-Null
    # Hyperlinks.Add re-styles the anchor cell with its own ad-hoc xf;
    # put it back on the sheet's shared "Hyperlink" cell style so it
    # matches every other linked cell above it.
    $ws.Range("B$row").Style = "Hyperlink"
}
